$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.679.54"
$ws.Cells.Item(2, 5).Value = "  -0.29%  "

$ws.Cells.Item(3, 4).Value = "2.731.40"
$ws.Cells.Item(3, 5).Value = "  -0.53%  "

$ws.Cells.Item(4, 5).Value = "  -0.04%  "

$ws.Cells.Item(5, 4).Value = "'561.98"
$ws.Cells.Item(5, 5).Value = "  -1.83%  "

$ws.Cells.Item(6, 4).Value = "'159.40"
$ws.Cells.Item(6, 5).Value = "  +1.49%  "

$ws.Cells.Item(7, 4).Value = "'1.00"

$ws.Cells.Item(8, 4).Value = "'0.597"
$ws.Cells.Item(8, 5).Value = "  -0.67%  "

$ws.Cells.Item(9, 5).Value = "  +0.21%  "

$ws.Cells.Item(10, 4).Value = "'0.166"
$ws.Cells.Item(10, 5).Value = "  +4.20%  "

$ws.Cells.Item(11, 5).Value = "  +3.37%  "

$ws.Cells.Item(12, 5).Value = "  -0.47%  "

$ws.Cells.Item(13, 4).Value = "3.213.50"
$ws.Cells.Item(13, 5).Value = "  -0.63%  "

$ws.Cells.Item(14, 4).Value = "'26.77"
$ws.Cells.Item(14, 5).Value = "  +1.43%  "

$ws.Cells.Item(15, 4).Value = "63.522.64"
$ws.Cells.Item(15, 5).Value = "  -0.07%  "

$ws.Cells.Item(16, 5).Value = "  -0.09%  "

$ws.Cells.Item(17, 4).Value = "2.735.47"
$ws.Cells.Item(17, 5).Value = "  -0.60%  "

$ws.Cells.Item(18, 4).Value = "'12.56"
$ws.Cells.Item(18, 5).Value = "  +3.33%  "

$ws.Cells.Item(19, 4).Value = "'4.73"
$ws.Cells.Item(19, 5).Value = "  -1.24%  "

$ws.Cells.Item(20, 4).Value = "'353.87"
$ws.Cells.Item(20, 5).Value = "  +0.03%  "

$ws.Cells.Item(21, 5).Value = "  -2.30%  "

$ws.Cells.Item(22, 5).Value = "  +0.25%  "

$ws.Cells.Item(23, 5).Value = "  -2.81%  "

$ws.Cells.Item(24, 4).Value = "'64.34"
$ws.Cells.Item(24, 5).Value = "  -1.23%  "

$ws.Cells.Item(25, 4).Value = "'0.169"
$ws.Cells.Item(25, 5).Value = "  +0.31%  "

$ws.Cells.Item(26, 4).Value = "'1.00"
$ws.Cells.Item(26, 5).Value = "  -0.04%  "

$ws.Cells.Item(27, 4).Value = "'8.36"
$ws.Cells.Item(27, 5).Value = "  -0.13%  "

$ws.Cells.Item(28, 4).Value = "0.0₃0901"
$ws.Cells.Item(28, 5).Value = "  +0.24%  "

$ws.Cells.Item(29, 4).Value = "'1.95"
$ws.Cells.Item(29, 5).Value = "  +1.08%  "

$ws.Cells.Item(30, 4).Value = "'7.16"
$ws.Cells.Item(30, 5).Value = "  +3.29%  "

$ws.Cells.Item(31, 5).Value = "  +10.83%  "

$ws.Cells.Item(32, 4).Value = "'165.72"
$ws.Cells.Item(32, 5).Value = "  -2.13%  "

$ws.Cells.Item(33, 2).Value = "EthereumClassic"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(33, 4).Value = "'20.01"
$ws.Cells.Item(33, 5).Value = "  -0.49%  "

$ws.Cells.Item(34, 2).Value = "NEARProtocol"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(34, 4).Value = "'4.88"
$ws.Cells.Item(34, 5).Value = "  +1.08%  "

$ws.Cells.Item(35, 4).Value = "'0.999"
$ws.Cells.Item(35, 5).Value = "  -0.04%  "

$ws.Cells.Item(36, 4).Value = "'1.46"
$ws.Cells.Item(36, 5).Value = "  +2.35%  "

$ws.Cells.Item(37, 4).Value = "'1.79"
$ws.Cells.Item(37, 5).Value = "  +0.84%  "

$ws.Cells.Item(38, 4).Value = "'0.971"
$ws.Cells.Item(38, 5).Value = "  -0.23%  "

$ws.Cells.Item(39, 4).Value = "'345.31"
$ws.Cells.Item(39, 5).Value = "  +6.15%  "

$ws.Cells.Item(40, 4).Value = "'6.26"
$ws.Cells.Item(40, 5).Value = "  +2.25%  "

$ws.Cells.Item(41, 4).Value = "'4.09"
$ws.Cells.Item(41, 5).Value = "  -0.81%  "

$ws.Cells.Item(42, 5).Value = "  -0.73%  "

$ws.Cells.Item(43, 4).Value = "'21.76"
$ws.Cells.Item(43, 5).Value = "  +2.48%  "

$ws.Cells.Item(44, 4).Value = "'21.03"
$ws.Cells.Item(44, 5).Value = "  -0.91%  "

$ws.Cells.Item(45, 4).Value = "'0.0582"
$ws.Cells.Item(45, 5).Value = "  -0.49%  "

$ws.Cells.Item(46, 4).Value = "'0.626"
$ws.Cells.Item(46, 5).Value = "  +0.63%  "

$ws.Cells.Item(47, 2).Value = "Stellar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(47, 4).Value = "'0.0999"
$ws.Cells.Item(47, 5).Value = "  -0.62%  "

$ws.Cells.Item(48, 2).Value = "VeChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(48, 4).Value = "'0.0249"
$ws.Cells.Item(48, 5).Value = "  -1.34%  "

$ws.Cells.Item(49, 5).Value = "  -0.04%  "

$ws.Cells.Item(50, 4).Value = "'131.76"
$ws.Cells.Item(50, 5).Value = "  -2.15%  "

$ws.Cells.Item(51, 4).Value = "'11.06"
$ws.Cells.Item(51, 5).Value = "  +0.25%  "
